# Update Unet values in the "cloud" sheet of confronto_modelli
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("cloud")

# Row 3 (Unet row): Tempo Training [min/epoch] 4.72 -> 6.1, Tempo inferenza [s] 55 -> 58
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "6.1"
$ws.Range("E3").Value = 58

# Restore active sheet / selection to match target (cloud tab selected, E4 active)
$ws.Range("E4").Select()

# "local" sheet selection changes from J5 to E9
$wsLocal = $wb.Worksheets.Item("local")
$wsLocal.Range("E9").Select()

$ws.Activate()
